$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-02 23:06:28"
$wsOverview.Columns.Item(5).ColumnWidth = 16.85
$wsOverview.Columns.Item(6).ColumnWidth = 16.85

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-02 23:06:24"
$wsZhCn.Columns.Item(3).ColumnWidth = 16.85

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-02 23:06:28"
$wsDeDe.Columns.Item(3).ColumnWidth = 16.85
